$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add pins for LCD.Keypad (4-bit variation) on rows 10-15 (D9-D14),
# mirroring the existing "Keypad" marker pattern in column B/C.
for ($r = 10; $r -le 15; $r++) {
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = "x"
    $cellB.Interior.Color = 65535

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = "LCD.Keypad"
}

# Rename the default cell style from the German "Stand." to "Normal"
# (as written by a non-German Excel build on re-save).
$standStyle = $wb.Styles.Item("Stand.")
$standStyle.Delete()
$wb.Styles.Add("Normal") | Out-Null

# Update the active selection to match the saved view state.
$ws.Range("G17").Select() | Out-Null
